$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = [double]"14.29746233333333"
$ws.Range("H2").Value = [double]"42.892387"
$ws.Range("I2").Value = [double]"0.1492463349041602"
$ws.Range("J2").Value = [double]"0.1492463349041601"
$ws.Range("M2").Value = [double]"0.701472"
$ws.Range("N2").Value = [double]"2.104416"
$ws.Range("O2").Value = [double]"0.001432026266286835"
$ws.Range("P2").Value = [double]"0.001432026266286835"
$ws.Range("Q2").Value = [double]"10.029269497888"
$ws.Range("R2").Value = [double]"90.263425480992"
$ws.Range("S2").Value = [double]"0.000213724671729799"
$ws.Range("T2").Value = [double]"0.000213724671729799"
$ws.Range("G3").Value = [double]"14.29746233333333"
$ws.Range("H3").Value = [double]"42.892387"
$ws.Range("I3").Value = [double]"0.1492463349041602"
$ws.Range("J3").Value = [double]"0.1492463349041601"
$ws.Range("O3").Value = [double]"0.004395913625032331"
$ws.Range("P3").Value = [double]"0.00439591362503233"
$ws.Range("Q3").Value = [double]"30.78700682579256"
$ws.Range("R3").Value = [double]"277.083061432133"
$ws.Range("S3").Value = [double]"0.000656073997091336"
$ws.Range("T3").Value = [double]"0.0006560739970913358"
$ws.Range("G4").Value = [double]"14.29746233333333"
$ws.Range("H4").Value = [double]"42.892387"
$ws.Range("I4").Value = [double]"0.1492463349041602"
$ws.Range("J4").Value = [double]"0.1492463349041601"
$ws.Range("M4").Value = [double]"243.6819663333333"
$ws.Range("N4").Value = [double]"731.045899"
$ws.Range("O4").Value = [double]"0.4974667219928344"
$ws.Range("P4").Value = [double]"0.4974667219928344"
$ws.Range("Q4").Value = [double]"3484.033734963435"
$ws.Range("R4").Value = [double]"31356.30361467091"
$ws.Range("S4").Value = [double]"0.07424508499421731"
$ws.Range("T4").Value = [double]"0.07424508499421728"
$ws.Range("G5").Value = [double]"14.29746233333333"
$ws.Range("H5").Value = [double]"42.892387"
$ws.Range("I5").Value = [double]"0.1492463349041602"
$ws.Range("J5").Value = [double]"0.1492463349041601"
$ws.Range("M5").Value = [double]"0.888026"
$ws.Range("N5").Value = [double]"2.664078"
$ws.Range("O5").Value = [double]"0.00181286859225405"
$ws.Range("P5").Value = [double]"0.00181286859225405"
$ws.Range("Q5").Value = [double]"12.69651828602067"
$ws.Range("R5").Value = [double]"114.268664574186"
$ws.Range("S5").Value = [double]"0.0002705639930567813"
$ws.Range("T5").Value = [double]"0.0002705639930567813"
$ws.Range("G6").Value = [double]"14.29746233333333"
$ws.Range("H6").Value = [double]"42.892387"
$ws.Range("I6").Value = [double]"0.1492463349041602"
$ws.Range("J6").Value = [double]"0.1492463349041601"
$ws.Range("M6").Value = [double]"242.4209796666667"
$ws.Range("N6").Value = [double]"727.262939"
$ws.Range("O6").Value = [double]"0.4948924695235924"
$ws.Range("P6").Value = [double]"0.4948924695235923"
$ws.Range("Q6").Value = [double]"3466.004825593932"
$ws.Range("R6").Value = [double]"31194.04343034539"
$ws.Range("S6").Value = [double]"0.07386088724806496"
$ws.Range("T6").Value = [double]"0.07386088724806493"
$ws.Range("I7").Value = [double]"0.7379913935691633"
$ws.Range("J7").Value = [double]"0.7379913935691632"
$ws.Range("M7").Value = [double]"0.701472"
$ws.Range("N7").Value = [double]"2.104416"
$ws.Range("O7").Value = [double]"0.001432026266286835"
$ws.Range("P7").Value = [double]"0.001432026266286835"
$ws.Range("Q7").Value = [double]"49.592605258816"
$ws.Range("R7").Value = [double]"446.333447329344"
$ws.Range("S7").Value = [double]"0.001056823059884667"
$ws.Range("T7").Value = [double]"0.001056823059884667"
$ws.Range("I8").Value = [double]"0.7379913935691633"
$ws.Range("J8").Value = [double]"0.7379913935691632"
$ws.Range("O8").Value = [double]"0.004395913625032331"
$ws.Range("P8").Value = [double]"0.00439591362503233"
$ws.Range("S8").Value = [double]"0.003244146422147282"
$ws.Range("T8").Value = [double]"0.003244146422147281"
$ws.Range("I9").Value = [double]"0.7379913935691633"
$ws.Range("J9").Value = [double]"0.7379913935691632"
$ws.Range("M9").Value = [double]"243.6819663333333"
$ws.Range("N9").Value = [double]"731.045899"
$ws.Range("O9").Value = [double]"0.4974667219928344"
$ws.Range("P9").Value = [double]"0.4974667219928344"
$ws.Range("Q9").Value = [double]"17227.80604936632"
$ws.Range("R9").Value = [double]"155050.2544442969"
$ws.Range("S9").Value = [double]"0.3671261594177754"
$ws.Range("T9").Value = [double]"0.3671261594177753"
$ws.Range("I10").Value = [double]"0.7379913935691633"
$ws.Range("J10").Value = [double]"0.7379913935691632"
$ws.Range("M10").Value = [double]"0.888026"
$ws.Range("N10").Value = [double]"2.664078"
$ws.Range("O10").Value = [double]"0.00181286859225405"
$ws.Range("P10").Value = [double]"0.00181286859225405"
$ws.Range("Q10").Value = [double]"62.78158340969467"
$ws.Range("R10").Value = [double]"565.034250687252"
$ws.Range("S10").Value = [double]"0.001337881418755334"
$ws.Range("T10").Value = [double]"0.001337881418755333"
$ws.Range("I11").Value = [double]"0.7379913935691633"
$ws.Range("J11").Value = [double]"0.7379913935691632"
$ws.Range("M11").Value = [double]"242.4209796666667"
$ws.Range("N11").Value = [double]"727.262939"
$ws.Range("O11").Value = [double]"0.4948924695235924"
$ws.Range("P11").Value = [double]"0.4948924695235923"
$ws.Range("Q11").Value = [double]"17138.6569258138"
$ws.Range("R11").Value = [double]"154247.9123323242"
$ws.Range("S11").Value = [double]"0.3652263832506006"
$ws.Range("T11").Value = [double]"0.3652263832506005"
$ws.Range("G12").Value = [double]"3.216815"
$ws.Range("H12").Value = [double]"9.650444999999999"
$ws.Range("I12").Value = [double]"0.03357923508533526"
$ws.Range("J12").Value = [double]"0.03357923508533525"
$ws.Range("M12").Value = [double]"0.701472"
$ws.Range("N12").Value = [double]"2.104416"
$ws.Range("O12").Value = [double]"0.001432026266286835"
$ws.Range("P12").Value = [double]"0.001432026266286835"
$ws.Range("Q12").Value = [double]"2.25650565168"
$ws.Range("R12").Value = [double]"20.30855086512"
$ws.Range("S12").Value = [double]"4.808634664402054E-05"
$ws.Range("T12").Value = [double]"4.808634664402053E-05"
$ws.Range("G13").Value = [double]"3.216815"
$ws.Range("H13").Value = [double]"9.650444999999999"
$ws.Range("I13").Value = [double]"0.03357923508533526"
$ws.Range("J13").Value = [double]"0.03357923508533525"
$ws.Range("O13").Value = [double]"0.004395913625032331"
$ws.Range("P13").Value = [double]"0.00439591362503233"
$ws.Range("Q13").Value = [double]"6.926831003528334"
$ws.Range("R13").Value = [double]"62.341479031755"
$ws.Range("S13").Value = [double]"0.000147611417029789"
$ws.Range("T13").Value = [double]"0.0001476114170297889"
$ws.Range("G14").Value = [double]"3.216815"
$ws.Range("H14").Value = [double]"9.650444999999999"
$ws.Range("I14").Value = [double]"0.03357923508533526"
$ws.Range("J14").Value = [double]"0.03357923508533525"
$ws.Range("M14").Value = [double]"243.6819663333333"
$ws.Range("N14").Value = [double]"731.045899"
$ws.Range("O14").Value = [double]"0.4974667219928344"
$ws.Range("P14").Value = [double]"0.4974667219928344"
$ws.Range("Q14").Value = [double]"783.8798045305616"
$ws.Range("R14").Value = [double]"7054.918240775054"
$ws.Range("S14").Value = [double]"0.01670455200492851"
$ws.Range("T14").Value = [double]"0.0167045520049285"
$ws.Range("G15").Value = [double]"3.216815"
$ws.Range("H15").Value = [double]"9.650444999999999"
$ws.Range("I15").Value = [double]"0.03357923508533526"
$ws.Range("J15").Value = [double]"0.03357923508533525"
$ws.Range("M15").Value = [double]"0.888026"
$ws.Range("N15").Value = [double]"2.664078"
$ws.Range("O15").Value = [double]"0.00181286859225405"
$ws.Range("P15").Value = [double]"0.00181286859225405"
$ws.Range("Q15").Value = [double]"2.85661535719"
$ws.Range("R15").Value = [double]"25.70953821471"
$ws.Range("S15").Value = [double]"6.087474063811953E-05"
$ws.Range("T15").Value = [double]"6.087474063811951E-05"
$ws.Range("G16").Value = [double]"3.216815"
$ws.Range("H16").Value = [double]"9.650444999999999"
$ws.Range("I16").Value = [double]"0.03357923508533526"
$ws.Range("J16").Value = [double]"0.03357923508533525"
$ws.Range("M16").Value = [double]"242.4209796666667"
$ws.Range("N16").Value = [double]"727.262939"
$ws.Range("O16").Value = [double]"0.4948924695235924"
$ws.Range("P16").Value = [double]"0.4948924695235923"
$ws.Range("Q16").Value = [double]"779.8234437064283"
$ws.Range("R16").Value = [double]"7018.410993357854"
$ws.Range("S16").Value = [double]"0.01661811057609482"
$ws.Range("T16").Value = [double]"0.01661811057609482"
$ws.Range("G17").Value = [double]"5.730442333333333"
$ws.Range("H17").Value = [double]"17.191327"
$ws.Range("I17").Value = [double]"0.059818133854125"
$ws.Range("J17").Value = [double]"0.059818133854125"
$ws.Range("M17").Value = [double]"0.701472"
$ws.Range("N17").Value = [double]"2.104416"
$ws.Range("O17").Value = [double]"0.001432026266286835"
$ws.Range("P17").Value = [double]"0.001432026266286835"
$ws.Range("Q17").Value = [double]"4.019744844448"
$ws.Range("R17").Value = [double]"36.177703600032"
$ws.Range("S17").Value = [double]"8.566113887936874E-05"
$ws.Range("T17").Value = [double]"8.566113887936873E-05"
$ws.Range("G18").Value = [double]"5.730442333333333"
$ws.Range("H18").Value = [double]"17.191327"
$ws.Range("I18").Value = [double]"0.059818133854125"
$ws.Range("J18").Value = [double]"0.059818133854125"
$ws.Range("O18").Value = [double]"0.004395913625032331"
$ws.Range("P18").Value = [double]"0.00439591362503233"
$ws.Range("Q18").Value = [double]"12.33947417506589"
$ws.Range("R18").Value = [double]"111.055267575593"
$ws.Range("S18").Value = [double]"0.0002629553496333558"
$ws.Range("T18").Value = [double]"0.0002629553496333557"
$ws.Range("G19").Value = [double]"5.730442333333333"
$ws.Range("H19").Value = [double]"17.191327"
$ws.Range("I19").Value = [double]"0.059818133854125"
$ws.Range("J19").Value = [double]"0.059818133854125"
$ws.Range("M19").Value = [double]"243.6819663333333"
$ws.Range("N19").Value = [double]"731.045899"
$ws.Range("O19").Value = [double]"0.4974667219928344"
$ws.Range("P19").Value = [double]"0.4974667219928344"
$ws.Range("Q19").Value = [double]"1396.405455746441"
$ws.Range("R19").Value = [double]"12567.64910171797"
$ws.Range("S19").Value = [double]"0.02975753096414016"
$ws.Range("T19").Value = [double]"0.02975753096414015"
$ws.Range("G20").Value = [double]"5.730442333333333"
$ws.Range("H20").Value = [double]"17.191327"
$ws.Range("I20").Value = [double]"0.059818133854125"
$ws.Range("J20").Value = [double]"0.059818133854125"
$ws.Range("M20").Value = [double]"0.888026"
$ws.Range("N20").Value = [double]"2.664078"
$ws.Range("O20").Value = [double]"0.00181286859225405"
$ws.Range("P20").Value = [double]"0.00181286859225405"
$ws.Range("Q20").Value = [double]"5.088781783500667"
$ws.Range("R20").Value = [double]"45.799036051506"
$ws.Range("S20").Value = [double]"0.0001084424161113919"
$ws.Range("T20").Value = [double]"0.0001084424161113919"
$ws.Range("G21").Value = [double]"5.730442333333333"
$ws.Range("H21").Value = [double]"17.191327"
$ws.Range("I21").Value = [double]"0.059818133854125"
$ws.Range("J21").Value = [double]"0.059818133854125"
$ws.Range("M21").Value = [double]"242.4209796666667"
$ws.Range("N21").Value = [double]"727.262939"
$ws.Range("O21").Value = [double]"0.4948924695235924"
$ws.Range("P21").Value = [double]"0.4948924695235923"
$ws.Range("Q21").Value = [double]"1389.179444370006"
$ws.Range("R21").Value = [double]"12502.61499933005"
$ws.Range("S21").Value = [double]"0.02960354398536073"
$ws.Range("T21").Value = [double]"0.02960354398536072"
$ws.Range("G22").Value = [double]"1.855114"
$ws.Range("H22").Value = [double]"5.565342"
$ws.Range("I22").Value = [double]"0.01936490258721644"
$ws.Range("J22").Value = [double]"0.01936490258721643"
$ws.Range("M22").Value = [double]"0.701472"
$ws.Range("N22").Value = [double]"2.104416"
$ws.Range("O22").Value = [double]"0.001432026266286835"
$ws.Range("P22").Value = [double]"0.001432026266286835"
$ws.Range("Q22").Value = [double]"1.301310527808"
$ws.Range("R22").Value = [double]"11.711794750272"
$ws.Range("S22").Value = [double]"2.773104914897982E-05"
$ws.Range("T22").Value = [double]"2.773104914897982E-05"
$ws.Range("G23").Value = [double]"1.855114"
$ws.Range("H23").Value = [double]"5.565342"
$ws.Range("I23").Value = [double]"0.01936490258721644"
$ws.Range("J23").Value = [double]"0.01936490258721643"
$ws.Range("O23").Value = [double]"0.004395913625032331"
$ws.Range("P23").Value = [double]"0.00439591362503233"
$ws.Range("Q23").Value = [double]"3.994653460108668"
$ws.Range("R23").Value = [double]"35.951881140978"
$ws.Range("S23").Value = [double]"8.512643913056856E-05"
$ws.Range("T23").Value = [double]"8.512643913056854E-05"
$ws.Range("G24").Value = [double]"1.855114"
$ws.Range("H24").Value = [double]"5.565342"
$ws.Range("I24").Value = [double]"0.01936490258721644"
$ws.Range("J24").Value = [double]"0.01936490258721643"
$ws.Range("M24").Value = [double]"243.6819663333333"
$ws.Range("N24").Value = [double]"731.045899"
$ws.Range("O24").Value = [double]"0.4974667219928344"
$ws.Range("P24").Value = [double]"0.4974667219928344"
$ws.Range("Q24").Value = [double]"452.0578272924953"
$ws.Range("R24").Value = [double]"4068.520445632458"
$ws.Range("S24").Value = [double]"0.009633394611773118"
$ws.Range("T24").Value = [double]"0.009633394611773115"
$ws.Range("G25").Value = [double]"1.855114"
$ws.Range("H25").Value = [double]"5.565342"
$ws.Range("I25").Value = [double]"0.01936490258721644"
$ws.Range("J25").Value = [double]"0.01936490258721643"
$ws.Range("M25").Value = [double]"0.888026"
$ws.Range("N25").Value = [double]"2.664078"
$ws.Range("O25").Value = [double]"0.00181286859225405"
$ws.Range("P25").Value = [double]"0.00181286859225405"
$ws.Range("Q25").Value = [double]"1.647389464964"
$ws.Range("R25").Value = [double]"14.826505184676"
$ws.Range("S25").Value = [double]"3.510602369242387E-05"
$ws.Range("T25").Value = [double]"3.510602369242386E-05"
$ws.Range("G26").Value = [double]"1.855114"
$ws.Range("H26").Value = [double]"5.565342"
$ws.Range("I26").Value = [double]"0.01936490258721644"
$ws.Range("J26").Value = [double]"0.01936490258721643"
$ws.Range("M26").Value = [double]"242.4209796666667"
$ws.Range("N26").Value = [double]"727.262939"
$ws.Range("O26").Value = [double]"0.4948924695235924"
$ws.Range("P26").Value = [double]"0.4948924695235923"
$ws.Range("Q26").Value = [double]"449.7185532733487"
$ws.Range("R26").Value = [double]"4047.466979460138"
$ws.Range("S26").Value = [double]"0.009583544463471346"
$ws.Range("T26").Value = [double]"0.009583544463471342"
